# Commit new script: Func_AssignGroupName
# Adds a new "BaseProduct" column (L) to the item-actions sheet with its
# sample value, and leaves the selection on K3 (matching the saved
# worksheet state after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column L.
$ws.Range("L1").Value = "BaseProduct"
$ws.Range("L2").Value = "U8LQ2E"

# Match the saved selection state.
$ws.Range("K3").Select() | Out-Null
